$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 466, pushing existing data (old rows 466-477) down to 469-480
$ws.Range("A466:R468").EntireRow.Insert()

# New data for the week of date 44448 (rows 466-468)
$ws.Range("A466").Value = 1
$ws.Range("B466").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C466").Value = "Arica y Parinacota"
$ws.Range("D466").Value = 44448
$ws.Range("E466").Value = 15
$ws.Range("F466").Value = 100112020
$ws.Range("G466").Value = "Tomate"
$ws.Range("H466").Value = "Larga vida"
$ws.Range("I466").Value = "Primera"
$ws.Range("J466").Value = 300
$ws.Range("K466").Value = 6000
$ws.Range("L466").Value = 7000
$ws.Range("M466").Value = 6500
$ws.Range("N466").Value = "$/caja 10 kilos"
$ws.Range("O466").Value = "Región de Arica y Parinacota"
$ws.Range("P466").Value = 650
$ws.Range("Q466").Value = 10
$ws.Range("R466").Value = "Hortaliza"

$ws.Range("A467").Value = 1
$ws.Range("B467").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C467").Value = "Arica y Parinacota"
$ws.Range("D467").Value = 44448
$ws.Range("E467").Value = 15
$ws.Range("F467").Value = 100112020
$ws.Range("G467").Value = "Tomate"
$ws.Range("H467").Value = "Larga vida"
$ws.Range("I467").Value = "Segunda"
$ws.Range("J467").Value = 350
$ws.Range("K467").Value = 5000
$ws.Range("L467").Value = 6000
$ws.Range("M467").Value = 5500
$ws.Range("N467").Value = "$/caja 10 kilos"
$ws.Range("O467").Value = "Región de Arica y Parinacota"
$ws.Range("P467").Value = 550
$ws.Range("Q467").Value = 10
$ws.Range("R467").Value = "Hortaliza"

$ws.Range("A468").Value = 1
$ws.Range("B468").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C468").Value = "Arica y Parinacota"
$ws.Range("D468").Value = 44448
$ws.Range("E468").Value = 15
$ws.Range("F468").Value = 100112020
$ws.Range("G468").Value = "Tomate"
$ws.Range("H468").Value = "Larga vida"
$ws.Range("I468").Value = "Tercera"
$ws.Range("J468").Value = 400
$ws.Range("K468").Value = 4000
$ws.Range("L468").Value = 5000
$ws.Range("M468").Value = 4500
$ws.Range("N468").Value = "$/caja 10 kilos"
$ws.Range("O468").Value = "Región de Arica y Parinacota"
$ws.Range("P468").Value = 450
$ws.Range("Q468").Value = 10
$ws.Range("R468").Value = "Hortaliza"
